# Refresh the crypto price / 1h-volume table with the latest scrape, and
# promote ShibaInu above InternetComputer(DFINITY) and EnergySwap above
# FraxShare to reflect their updated ranking (rows 20/21 and 46/47).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell address -> new value, in document order.
$cellUpdates = [ordered]@{
    "D2" = "43.034.35"
    "E2" = "  +0.90%  "
    "D3" = "2.299.96"
    "E3" = "  +0.57%  "
    "D4" = "1.00"
    "E4" = "  +0.00%  "
    "D5" = "299.73"
    "E5" = "  -0.28%  "
    "D6" = "97.64"
    "E6" = "  +0.66%  "
    "E7" = "  +0.46%  "
    "E9" = "  +1.09%  "
    "D10" = "33.77"
    "E10" = "  +1.45%  "
    "D11" = "0.0793"
    "E11" = "  +0.61%  "
    "D12" = "49.03"
    "E12" = "  -1.99%  "
    "E13" = "  +2.78%  "
    "E14" = "  +10.48%  "
    "D15" = "6.78"
    "E15" = "  +1.88%  "
    "D16" = "2.658.22"
    "E16" = "  +0.63%  "
    "D17" = "2.299.57"
    "E17" = "  +1.06%  "
    "E18" = "  +2.97%  "
    "D19" = "42.964.13"
    "E19" = "  +0.93%  "
    "B20" = "ShibaInu"
    "C20" = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
    "D20" = "0.0₃0903"
    "E20" = "  +0.86%  "
    "B21" = "InternetComputer(DFINITY)"
    "C21" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D21" = "11.62"
    "E21" = "  +1.02%  "
    "E22" = "  +0.54%  "
    "E23" = "  +1.01%  "
    "D24" = "236.68"
    "E24" = "  +0.90%  "
    "D25" = "2.04"
    "E25" = "  +4.97%  "
    "E26" = "  +0.03%  "
    "D27" = "2.46"
    "E27" = "  -1.35%  "
    "E28" = "  -0.30%  "
    "D29" = "166.27"
    "E29" = "  +0.20%  "
    "E30" = "  +0.29%  "
    "E31" = "  +0.03%  "
    "D32" = "9.13"
    "E32" = "  +0.22%  "
    "E33" = "  +0.05%  "
    "D34" = "4.97"
    "E34" = "  -0.02%  "
    "E35" = "  +5.78%  "
    "E36" = "  +1.74%  "
    "D37" = "16.83"
    "E37" = "  +3.79%  "
    "D38" = "0.0703"
    "E38" = "  +1.09%  "
    "D39" = "2.83"
    "E39" = "  -0.20%  "
    "E40" = "  +0.80%  "
    "E41" = "  +0.80%  "
    "E42" = "  -0.37%  "
    "E43" = "  -4.37%  "
    "D44" = "1.993.98"
    "E44" = "  +1.62%  "
    "E45" = "  +0.67%  "
    "B46" = "EnergySwap"
    "C46" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D46" = "17.61"
    "E46" = "  -1.55%  "
    "B47" = "FraxShare"
    "C47" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D47" = "9.81"
    "E47" = "  +1.08%  "
    "D48" = "2.85"
    "E48" = "  +0.82%  "
    "D49" = "2.524.09"
    "E49" = "  +0.65%  "
    "D50" = "53.36"
    "E50" = "  +0.35%  "
    "D51" = "4.58"
    "E51" = "  -1.89%  "
}

foreach ($addr in $cellUpdates.Keys) {
    $newValue = $cellUpdates[$addr]
    $range = $ws.Range($addr)

    if ($newValue -match "^-?[0-9]+(\.[0-9]+)?$") {
        # This text looks like a plain number ("1.00", "33.77", ...). Force the
        # cell to Text format before assigning so Excel keeps it as a string
        # (preserving formatting such as trailing zeros) instead of silently
        # converting it to a numeric value. Reset the style afterwards so no
        # extra formatting is left behind on the cell.
        $range.NumberFormat = "@"
        $range.Value = $newValue
        $range.Style = "Normal"
    } else {
        $range.Value = $newValue
    }
}
